# UniDebAutomatedCar CommunicationMatrix.xlsx
# Add a new "Version & History" row (V1.6) documenting the change to the
# "Type" and "Receiver" columns for the "Vehicle Speed" and "Motor RPM"
# signals, and restore the "Version & History" tab as the active sheet.

$wb = $excel.ActiveWorkbook
$wsHist   = $wb.Worksheets.Item("Version & History")
$wsMatrix = $wb.Worksheets.Item("CommunicationMatrix")

# --- Append the new history row (row 11), cloning row 10's formatting ---
$wsHist.Range("A10:E10").Copy()
$wsHist.Range("A11:E11").PasteSpecial(-4122)

$wsHist.Range("A11").Value = "V1.6"
$wsHist.Range("B11").Value = "Changed the „Type” and „Receiver” columns for the „Vehicle Speed” and „Motor RPM” signals."
$wsHist.Range("C11").Value = "Zborai Attila"
$wsHist.Range("D11").Value = 42804
$wsHist.Range("E11").Value = "Draft version"

# --- Row heights (auto-recalculated by Excel when the text reflows) ---
$wsHist.Rows(7).RowHeight = 85.45
$wsHist.Rows(10).RowHeight = 20.85
$wsHist.Rows(11).RowHeight = 20.85

# --- Restore selections / active sheet to match the saved workbook state ---
$wsMatrix.Activate()
$wsMatrix.Range("A2").Select()

$wsHist.Activate()
$wsHist.Range("C14").Select()
